# Auto-generated Excel COM-interop script
# Applies updated "想去人数" (F column) values across 4 worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 54
$ws.Range("F6").Value = 2758
$ws.Range("F8").Value = 1644
$ws.Range("F9").Value = 7463
$ws.Range("F11").Value = 7659
$ws.Range("F15").Value = 6186
$ws.Range("F16").Value = 3267
$ws.Range("F17").Value = 3632
$ws.Range("F18").Value = 21
$ws.Range("F19").Value = 13
$ws.Range("F24").Value = 286
$ws.Range("F26").Value = 3639
$ws.Range("F28").Value = 342
$ws.Range("F29").Value = 927
$ws.Range("F30").Value = 262
$ws.Range("F31").Value = 1101
$ws.Range("F32").Value = 65
$ws.Range("F34").Value = 2622
$ws.Range("F35").Value = 1471
$ws.Range("F37").Value = 20
$ws.Range("F38").Value = 29
$ws.Range("F39").Value = 3292
$ws.Range("F40").Value = 168
$ws.Range("F41").Value = 245
$ws.Range("F43").Value = 901
$ws.Range("F45").Value = 1289
$ws.Range("F46").Value = 227
$ws.Range("F47").Value = 525
$ws.Range("F48").Value = 595

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 57
$ws.Range("F10").Value = 31
$ws.Range("F16").Value = 81

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 121

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 57
$ws.Range("F7").Value = 54
$ws.Range("F8").Value = 121
$ws.Range("F9").Value = 1644
$ws.Range("F12").Value = 7463
$ws.Range("F13").Value = 7659
$ws.Range("F15").Value = 6186
$ws.Range("F16").Value = 3267
$ws.Range("F17").Value = 3632
$ws.Range("F18").Value = 21
$ws.Range("F19").Value = 13
$ws.Range("F23").Value = 286
$ws.Range("F27").Value = 3639
$ws.Range("F30").Value = 342
$ws.Range("F31").Value = 927
$ws.Range("F32").Value = 262
$ws.Range("F33").Value = 65
$ws.Range("F35").Value = 2622
$ws.Range("F36").Value = 1471
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 81
$ws.Range("F40").Value = 3292
$ws.Range("F41").Value = 168
$ws.Range("F42").Value = 245
$ws.Range("F45").Value = 901
$ws.Range("F47").Value = 1289
$ws.Range("F48").Value = 227
$ws.Range("F49").Value = 525
